$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.695.16"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "'1.638.88"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'212.70"
$ws.Range("D6").Value = "'0.524"
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'23.17"
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("D9").Value = "'0.260"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "'1.870.11"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "'1.642.70"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("E15").Value = "  -4.33%  "
$ws.Range("D16").Value = "'64.74"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "'27.659.17"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "'230.51"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.31"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "'10.26"
$ws.Range("E23").Value = "  +4.88%  "
$ws.Range("D24").Value = "'2.07"
$ws.Range("E24").Value = "  +2.91%  "
$ws.Range("D25").Value = "'150.75"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "'15.60"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "'1.459.34"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "'0.566"
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "'0.895"
$ws.Range("E40").Value = "  +9.42%  "
$ws.Range("D41").Value = "'69.29"
$ws.Range("E41").Value = "  +5.92%  "
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").Value = "'5.58"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").Value = "'1.780.38"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("D49").Value = "'86.79"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0107"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0993"
$ws.Range("E51").Value = "  +0.07%  "
